# Apply cryptos list refresh (prices + 1h volume deltas) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.830.66'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '1.855.37'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '304.42'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = "'0.5040"
$ws.Range('E7').Value = '  -2.52%  '
$ws.Range('D8').Value = '0.3651'
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').Value = '0.8896'
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').Value = '20.65'
$ws.Range('E11').Value = '  -0.17%  '
$ws.Range('D12').Value = '0.07518'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '1.856.59'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = "'92.00"
$ws.Range('E14').Value = '  +2.89%  '
$ws.Range('E15').Value = '  -1.97%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').Value = "'0.000008500"
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '26.865.80'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').Value = '5.025'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').Value = '2.085.39'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('D23').Value = '10.32'
$ws.Range('E23').Value = '  -2.71%  '
$ws.Range('D24').Value = '6.449'
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').Value = '146.43'
$ws.Range('E25').Value = '  -3.00%  '
$ws.Range('D26').Value = '1.796'
$ws.Range('E26').Value = '  -2.75%  '
$ws.Range('D27').Value = '17.82'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '2.052'
$ws.Range('E28').Value = '  -4.80%  '
$ws.Range('D29').Value = '112.89'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '4.632'
$ws.Range('E30').Value = '  -2.36%  '
$ws.Range('D31').Value = '4.657'
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '0.09198'
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('D33').Value = '0.05089'
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D34').Value = '2.994'
$ws.Range('E34').Value = '  -3.46%  '
$ws.Range('D35').Value = '0.7341'
$ws.Range('E35').Value = '  -2.76%  '
$ws.Range('D36').Value = '1.143'
$ws.Range('E36').Value = '  -2.62%  '
$ws.Range('D37').Value = '3.233'
$ws.Range('E37').Value = '  +6.93%  '
$ws.Range('D38').Value = '0.01991'
$ws.Range('E38').Value = '  -2.03%  '
$ws.Range('D39').Value = '2.494'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').Value = '0.5314'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('D42').Value = '118.92'
$ws.Range('E42').Value = '  +3.10%  '
$ws.Range('D43').Value = '6.476'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('D44').Value = '8.352'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').Value = '0.1466'
$ws.Range('E45').Value = '  -1.22%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = "'1.000"
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.4631'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').Value = '9.932'
$ws.Range('E48').Value = '  -2.00%  '
$ws.Range('D49').Value = '1.555'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').Value = '62.81'
